$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting Code/Description/Definition right.
$ws.Columns("A").Insert()

# New header for the inserted column.
$ws.Range("A1").Value = "Version"

# Definition header moves to D1 (was blank before, now a header is needed there too
# since the old C1 "Definition" header shifted to D1 automatically via the insert).

# Fill the new "Version" column with "1.0" for every data row.
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = "1.0"
}

# sheetFormatPr gains baseColWidth="10".
$ws.StandardWidth = 10
